$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Row = 3;  F = 8569; G = $null },
    @{ Row = 4;  F = 1526; G = $null },
    @{ Row = 7;  F = 266;  G = 60 },
    @{ Row = 8;  F = 170;  G = $null },
    @{ Row = 9;  F = 34;   G = $null },
    @{ Row = 13; F = 1269; G = $null },
    @{ Row = 14; F = 281;  G = $null },
    @{ Row = 15; F = 84;   G = $null },
    @{ Row = 16; F = 149;  G = $null },
    @{ Row = 17; F = 103;  G = $null },
    @{ Row = 18; F = 137;  G = $null },
    @{ Row = 19; F = 83;   G = $null },
    @{ Row = 20; F = 125;  G = $null },
    @{ Row = 21; F = 112;  G = $null }
)

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($u in $updates) {
        $ws.Cells.Item($u.Row, 6).Value = $u.F
        if ($null -ne $u.G) {
            $ws.Cells.Item($u.Row, 7).Value = $u.G
        }
    }
}
